$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Servers")

# Remove the stale ping/comment notes (their shared strings get garbage
# collected on save, which also re-indexes the remaining Core/RAM(G)/Ping
# strings referenced from row 1).
$ws.Range("H4").ClearContents()
$ws.Range("I3").ClearContents()
$ws.Range("I4").ClearContents()

# Row 4 no longer needs the tall wrapped-text row.
$ws.Rows.Item(4).AutoFit()

# New "TW DEMO" highlight cell (yellow fill) + matching fill on the rest of
# row 3.
$ws.Range("A3:F3").Interior.Color = 65535
$ws.Range("F3").Value = "TW DEMO"

# New "For use by JU" note cell.
$ws.Range("F4").Value = "For use by JU: 27-29-Feg"

# Highlight the 18.222.6.115 server name with the accent font color.
$ws.Range("B4").Font.ThemeColor = 6
$ws.Range("B4").Font.TintAndShade = -0.249977111117893

# Selection moved from B7 to C4 on the Servers sheet.
$ws.Range("C4").Select()
